# Chain Azure Functions together using input and output bindings
# Mark row 32 ("Stage a web app deployment...") and row 34 ("Expose multiple
# Azure Function apps...") as Completed, with start/completed dates and
# Badge/Github hyperlinks - matching the formatting already used by the
# other "Completed" rows (e.g. row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the cell formatting (fill/border/font) used by an existing
#     "Completed" row onto the B:G cells of rows 32 and 34 -----------------
$ws.Range("B30:G30").Copy()
$ws.Range("B32:G32").PasteSpecial(-4122)
$ws.Range("B30:G30").Copy()
$ws.Range("B34:G34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 32: "Stage a web app deployment for testing and rollback by
#              using App Service deployment slot" ---------------------------
$ws.Range("C32").Value = "Completed"
$ws.Range("D32").Value = 45341
$ws.Range("E32").Value = 45341

# Shared-string table order: Github link (G32) before Badge link (F32).
$ws.Range("G32").Value = "https://github.com/ShubhamVermaTheDeveloper/AzureModuleLearning/tree/main/Stage%20a%20web%20app%20deployment%20for%20testing%20and%20rollback%20by%20using%20App%20Service%20deployment%20slot"
$ws.Range("F32").Value = "https://learn.microsoft.com/api/achievements/share/en-us/ShubhamVerma/AQW9HSD7?sharingId=7BBBB75FB7AF740D"

$ws.Hyperlinks.Add($ws.Range("G32"), "https://github.com/ShubhamVermaTheDeveloper/AzureModuleLearning/tree/main/Stage%20a%20web%20app%20deployment%20for%20testing%20and%20rollback%20by%20using%20App%20Service%20deployment%20slot", "", "", "")
$ws.Hyperlinks.Add($ws.Range("F32"), "https://learn.microsoft.com/api/achievements/share/en-us/ShubhamVerma/AQW9HSD7?sharingId=7BBBB75FB7AF740D", "", "", "")

# Restore the Badge-link/Github-link formatting that Hyperlinks.Add resets.
$ws.Range("G30").Copy()
$ws.Range("G32").PasteSpecial(-4122)
$ws.Range("F30").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 34: "Expose multiple Azure Function apps as a consistent API by
#              using Azure API Management" -----------------------------------
$ws.Range("C34").Value = "Completed"
$ws.Range("D34").Value = 45341
$ws.Range("E34").Value = 45341

# Shared-string table order: Badge link (F34) before Github link (G34).
$ws.Range("F34").Value = "https://learn.microsoft.com/api/achievements/share/en-us/ShubhamVerma/PTZYCFM4?sharingId=7BBBB75FB7AF740D"
$ws.Range("G34").Value = "https://github.com/ShubhamVermaTheDeveloper/AzureModuleLearning/tree/main/Expose%20multiple%20Azure%20Function%20apps%20as%20a%20consistent%20API%20by%20using%20Azure%20API%20Management"

$ws.Hyperlinks.Add($ws.Range("F34"), "https://learn.microsoft.com/api/achievements/share/en-us/ShubhamVerma/PTZYCFM4?sharingId=7BBBB75FB7AF740D", "", "", "")
$ws.Hyperlinks.Add($ws.Range("G34"), "https://github.com/ShubhamVermaTheDeveloper/AzureModuleLearning/tree/main/Expose%20multiple%20Azure%20Function%20apps%20as%20a%20consistent%20API%20by%20using%20Azure%20API%20Management", "", "", "")

# Restore the Badge-link/Github-link formatting that Hyperlinks.Add resets.
$ws.Range("F30").Copy()
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("G30").Copy()
$ws.Range("G34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the saved view state (scroll position / zoom / selection) ------
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Application.ActiveWindow.Zoom = 93
$ws.Range("E47").Select()
